{"js": "// Change Daniel Engbert's resume contact email from his personal Gmail\n// address to his UMBC school address (commit: \"changed email to UMBC email\").\n\nconst body = context.document.body;\n\n// Find the old email address in the document body.\nconst results = body.search(\"danielengbert@gmail.com\", {\n  matchCase: false,\n  matchWholeWord: false\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"end1@umbc.edu\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Change Daniel Engbert's resume contact email from his personal Gmail\n# address to his UMBC school address (commit: \"changed email to UMBC email\").\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace\n$range.Find.Execute(\"danielengbert@gmail.com\", $false, $false, $false, $false, $false, $true, 1, $false, \"end1@umbc.edu\", 2)\n"}
